# Refresh the cryptos list (Price / Volume(1h) columns) with the latest
# scraped values. Price cells (column D) are forced to Text format before
# the assignment so values such as "26.914.49" / "0.07191" are stored
# verbatim instead of being auto-parsed into numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.914.49'
$ws.Range("E2").Value = '  -1.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.831.93'
$ws.Range("E3").Value = '  -2.19%  '
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.66'
$ws.Range("E5").Value = '  -1.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("E7").Value = '  -1.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3654'
$ws.Range("E8").Value = '  -2.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07191'
$ws.Range("E9").Value = '  -2.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8791'
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07870'
$ws.Range("E11").Value = '  -1.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.66'
$ws.Range("E12").Value = '  -2.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.914.70'
$ws.Range("E13").Value = '  +3.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.321'
$ws.Range("E14").Value = '  -2.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.422'
$ws.Range("E15").Value = '  -3.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.71'
$ws.Range("E16").Value = '  -3.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.007'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008741'
$ws.Range("E18").Value = '  -2.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.933.03'
$ws.Range("E20").Value = '  -2.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.51'
$ws.Range("E21").Value = '  -3.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.005'
$ws.Range("E22").Value = '  -3.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.42'
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.007'
$ws.Range("E24").Value = '  +6.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.42'
$ws.Range("E25").Value = '  -1.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.18'
$ws.Range("E26").Value = '  -2.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.988'
$ws.Range("E27").Value = '  -6.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '114.26'
$ws.Range("E28").Value = '  -2.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.939'
$ws.Range("E29").Value = '  -4.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08812'
$ws.Range("E30").Value = '  -1.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.099'
$ws.Range("E31").Value = '  +3.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7634'
$ws.Range("E32").Value = '  +0.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.464'
$ws.Range("E33").Value = '  -1.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.134'
$ws.Range("E34").Value = '  -2.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.609'
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.093'
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01927'
$ws.Range("E37").Value = '  -2.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05156'
$ws.Range("E38").Value = '  -2.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.928'
$ws.Range("E39").Value = '  -2.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.940'
$ws.Range("E40").Value = '  -4.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4988'
$ws.Range("E41").Value = '  -5.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1601'
$ws.Range("E42").Value = '  -3.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.367'
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.30'
$ws.Range("E44").Value = '  -1.29%  '
$ws.Range("E45").Value = '  -5.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.006'
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.59'
$ws.Range("E47").Value = '  -1.50%  '
$ws.Range("E48").Value = '  -4.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06112'
$ws.Range("E49").Value = '  -2.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '64.68'
$ws.Range("E50").Value = '  -2.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.30'
$ws.Range("E51").Value = '  -2.85%  '
